# brand crawling 공백 제거
# Populate the (previously empty) sheet with the crawled brand/artist rows.
# Numeric-looking values that must stay as TEXT are written with a leading
# apostrophe so Excel keeps them as strings instead of auto-converting them
# to numbers (matches the source data, which stores them as inline strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 1
$ws.Cells.Item(1, 2).Value = "'"
$ws.Cells.Item(1, 3).Value = 7
$ws.Cells.Item(1, 4).Value = "아티스트명"
$ws.Cells.Item(1, 5).Value = "A Tiseuteumyung"
$ws.Cells.Item(1, 6).Value = "'3114174"

# Row 2
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 7
$ws.Cells.Item(2, 4).Value = "민지"
$ws.Cells.Item(2, 5).Value = "'3116664"

# Row 3
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 7
$ws.Cells.Item(3, 4).Value = "하니"
$ws.Cells.Item(3, 5).Value = "'3116665"
